$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.013.62"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.615.81"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'594.37"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'166.56"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").Value = "2.616.40"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "'0.362"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "'5.22"
$ws.Range("D14").Value = "'27.61"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "3.103.81"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "'0.0000182"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "67.294.80"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "2.615.42"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'11.98"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").Value = "'7.98"
$ws.Range("E20").Value = "  +6.96%  "
$ws.Range("D21").Value = "'356.02"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'10.25"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("E26").Value = "  -4.47%  "
$ws.Range("D27").Value = "'69.61"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "0.0₃0998"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "'546.60"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'7.87"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").Value = "'1.89"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  +6.31%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "'158.10"
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "'0.365"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "'18.14"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "'5.18"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'2.40"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").Value = "0.0₆0296"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").Value = "'151.38"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").Value = "'1.68"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("E51").Value = "  -0.65%  "
